$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.717.53"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.91%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.963.63"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.26%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.77"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.61%  "

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.85%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.17"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.57%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.09%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0818"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.99%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.22%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.39"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +4.41%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.253.78"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.28%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.829"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.05%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "13.74"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.05%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.27"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.75%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.962.98"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.20%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.597.54"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.70%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.29%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0860"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.35%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.08"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.65%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.63"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.23%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.02%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.37%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +3.00%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +16.05%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.30"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.64%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "160.77"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.74%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.37"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.02%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.65%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.54%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.28%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0621"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.85%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.30"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.79%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +6.82%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.20%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.06"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.37%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.40"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +12.05%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.27%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0988"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.35%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.05%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.01%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.71%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.21"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.96%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.361.36"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.75%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.98%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "87.85"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.38%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.17"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.46%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.87%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.144.16"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.26%  "

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.64%  "
